$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force these price cells to Text format before assignment so Excel
# does not silently reinterpret the scraped price strings (which look
# like plain decimals) as numeric values and normalize/round them.
$ws.Range("D5,D6,D8,D10,D11,D12,D13,D14,D16,D17,D21,D22,D23,D24,D25,D26,D27,D29,D30,D31,D32,D34,D35,D36,D37,D39,D40,D41,D43,D44,D46,D47,D48,D50,D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "70.655.48"
$ws.Range("E2").Value = "  -1.91%  "

# Row 3
$ws.Range("D3").Value = "3.636.47"
$ws.Range("E3").Value = "  +1.03%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "581.30"
$ws.Range("E5").Value = "  -2.00%  "

# Row 6
$ws.Range("D6").Value = "175.66"
$ws.Range("E6").Value = "  -4.54%  "

# Row 7
$ws.Range("D7").Value = "3.627.64"
$ws.Range("E7").Value = "  +0.92%  "

# Row 8
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  +0.21%  "

# Row 9
$ws.Range("E9").Value = "  +0.15%  "

# Row 10
$ws.Range("D10").Value = "0.197"
$ws.Range("E10").Value = "  -5.08%  "

# Row 11
$ws.Range("D11").Value = "6.88"
$ws.Range("E11").Value = "  +19.79%  "

# Row 12
$ws.Range("D12").Value = "0.605"
$ws.Range("E12").Value = "  -0.60%  "

# Row 13
$ws.Range("D13").Value = "48.42"
$ws.Range("E13").Value = "  -3.63%  "

# Row 14
$ws.Range("D14").Value = "0.0000285"
$ws.Range("E14").Value = "  -1.85%  "

# Row 15
$ws.Range("D15").Value = "4.222.33"
$ws.Range("E15").Value = "  +1.34%  "

# Row 16
$ws.Range("D16").Value = "666.61"
$ws.Range("E16").Value = "  -4.49%  "

# Row 17
$ws.Range("D17").Value = "8.90"
$ws.Range("E17").Value = "  -0.11%  "

# Row 18
$ws.Range("D18").Value = "3.625.49"
$ws.Range("E18").Value = "  +1.41%  "

# Row 19
$ws.Range("D19").Value = "70.637.96"
$ws.Range("E19").Value = "  -1.94%  "

# Row 20
$ws.Range("E20").Value = "  -0.14%  "

# Row 21
$ws.Range("D21").Value = "17.78"
$ws.Range("E21").Value = "  -2.89%  "

# Row 22
$ws.Range("D22").Value = "11.40"
$ws.Range("E22").Value = "  -3.21%  "

# Row 23
$ws.Range("D23").Value = "0.939"
$ws.Range("E23").Value = "  +0.93%  "

# Row 24
$ws.Range("D24").Value = "17.12"
$ws.Range("E24").Value = "  -3.50%  "

# Row 25
$ws.Range("D25").Value = "99.72"
$ws.Range("E25").Value = "  -4.82%  "

# Row 26
$ws.Range("D26").Value = "3.92"
$ws.Range("E26").Value = "  -2.50%  "

# Row 27
$ws.Range("D27").Value = "2.79"
$ws.Range("E27").Value = "  -1.06%  "

# Row 28
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -0.72%  "

# Row 30
$ws.Range("D30").Value = "34.71"
$ws.Range("E30").Value = "  -1.89%  "

# Row 31
$ws.Range("D31").Value = "3.35"
$ws.Range("E31").Value = "  -3.75%  "

# Row 32
$ws.Range("D32").Value = "8.99"
$ws.Range("E32").Value = "  -0.83%  "

# Row 33
$ws.Range("E33").Value = "  -5.83%  "

# Row 34
$ws.Range("D34").Value = "7.40"
$ws.Range("E34").Value = "  -1.21%  "

# Row 35
$ws.Range("D35").Value = "4.00"
$ws.Range("E35").Value = "  -3.44%  "

# Row 36
$ws.Range("D36").Value = "587.42"
$ws.Range("E36").Value = "  -0.32%  "

# Row 37
$ws.Range("D37").Value = "11.07"
$ws.Range("E37").Value = "  -2.47%  "

# Row 38
$ws.Range("E38").Value = "  +0.08%  "

# Row 39
$ws.Range("D39").Value = "58.21"
$ws.Range("E39").Value = "  -2.69%  "

# Row 40
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  +0.02%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0456"
$ws.Range("E41").Value = "  +3.20%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.572.21"
$ws.Range("E42").Value = "  -2.85%  "

# Row 43
$ws.Range("D43").Value = "0.141"
$ws.Range("E43").Value = "  -3.21%  "

# Row 44
$ws.Range("D44").Value = "0.344"
$ws.Range("E44").Value = "  -1.24%  "

# Row 45
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0739"
$ws.Range("E45").Value = "  -5.93%  "

# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "34.50"
$ws.Range("E46").Value = "  -5.03%  "

# Row 47
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").Value = "  -4.61%  "

# Row 48
$ws.Range("D48").Value = "2.92"
$ws.Range("E48").Value = "  +5.14%  "

# Row 49
$ws.Range("E49").Value = "  +0.35%  "

# Row 50
$ws.Range("D50").Value = "135.66"
$ws.Range("E50").Value = "  +0.91%  "

# Row 51
$ws.Range("D51").Value = "2.97"
$ws.Range("E51").Value = "  +2.22%  "

# Restore the cell style/format marker to the workbook default ("Normal")
# now that the text values are safely stored, so these cells end up with
# no explicit style index, matching the rest of the sheet.
$ws.Range("D5,D6,D8,D10,D11,D12,D13,D14,D16,D17,D21,D22,D23,D24,D25,D26,D27,D29,D30,D31,D32,D34,D35,D36,D37,D39,D40,D41,D43,D44,D46,D47,D48,D50,D51").Style = "Normal"
